# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values on the zh-cn and
# de-de report sheets to reflect a newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-11 02:44:19"
$wsZhCn.Range("G2").Value = "2016-01-11 02:45:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-11 02:44:33"
$wsDeDe.Range("G2").Value = "2016-01-11 02:45:32"
